$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76 currently holds one week's price record. The edit both updates
# row 76 to a newer week's figures AND inserts the old figures (that used
# to live in row 76) as a brand-new row 77, so the sheet grows by one row.

# 1) Snapshot row 76's current (pre-edit) values across all 18 columns
#    (A..R) using Value2 (a plain, non-parameterized property so it reads
#    cleanly through COM interop).
$orig = @{}
for ($col = 1; $col -le 18; $col++) {
    $orig[$col] = $ws.Cells.Item(76, $col).Value2
}

# 2) Write that snapshot into the new row 77, preserving values as-is.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(77, $col).Value2 = $orig[$col]
}
# Match the date cell's number formatting to the source cell (D76) so the
# new D77 date renders the same way.
$ws.Cells.Item(77, 4).NumberFormat = $ws.Cells.Item(76, 4).NumberFormat

# 3) Update row 76 in place with the new week's values.
$ws.Cells.Item(76, 4).Value2 = 44656    # D76 Fecha
$ws.Cells.Item(76, 11).Value2 = 8500    # K76 Precio minimo
$ws.Cells.Item(76, 12).Value2 = 9000    # L76 Precio maximo
$ws.Cells.Item(76, 13).Value2 = 8750    # M76 Precio promedio ponderado
$ws.Cells.Item(76, 16).Value2 = 583     # P76 Precio $/Kg
